# New words and +2Lvl`s
# Adds two new entry rows (11 & 12) to the word list on the active sheet
# and moves the active selection, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Медіапростір"
$ws.Range("B11").Value = "entry.1568295105"

$ws.Range("A12").Value = "Місцеперебування"
$ws.Range("B12").Value = "entry.283921628"

$ws.Range("H16").Select() | Out-Null
